$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")
$ws.Activate()

# Update End Year value from 2025 to 2024
$ws.Range("B4").Value = 2024

# Update maximum_investment_capacity_per_year value from 1000 to 4000
$ws.Range("B13").Value = 4000

# Update the selected cell on the sheet to C3
$ws.Range("C3").Select()
